# Updated cryptos list with latest price/volume data from coinranking.com
# Two rows (WEMIXTOKEN/FraxShare and Flow/BabyDogeCoin) also swapped rank position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ Col = NewValue; ... } for every cell that changed in this update
$rowUpdates = [ordered]@{
    2 = @{ "D" = "24.625.33"; "E" = "  +1.60%  " }
    3 = @{ "D" = "1.696.22"; "E" = "  +2.01%  " }
    4 = @{ "D" = "0.9993"; "E" = "  -0.77%  " }
    5 = @{ "D" = "313.00"; "E" = "  +0.74%  " }
    6 = @{ "D" = "0.9986"; "E" = "  -0.80%  " }
    7 = @{ "D" = "0.3948"; "E" = "  +0.98%  " }
    8 = @{ "D" = "0.4033"; "E" = "  +2.28%  " }
    9 = @{ "D" = "0.9991"; "E" = "  -0.79%  " }
    10 = @{ "D" = "1.519"; "E" = "  +8.97%  " }
    11 = @{ "D" = "54.34"; "E" = "  +11.64%  " }
    12 = @{ "D" = "0.08753"; "E" = "  +1.71%  " }
    13 = @{ "D" = "7.322"; "E" = "  +13.67%  " }
    14 = @{ "D" = "23.12"; "E" = "  +2.14%  " }
    15 = @{ "D" = "0.00001319"; "E" = "  +2.51%  " }
    16 = @{ "D" = "7.596"; "E" = "  +6.99%  " }
    17 = @{ "D" = "1.695.27"; "E" = "  +2.06%  " }
    18 = @{ "D" = "100.27"; "E" = "  -0.70%  " }
    19 = @{ "D" = "0.07076"; "E" = "  +4.61%  " }
    20 = @{ "D" = "19.45"; "E" = "  +3.22%  " }
    21 = @{ "D" = "6.710"; "E" = "  +1.57%  " }
    22 = @{ "D" = "0.9998"; "E" = "  -0.61%  " }
    23 = @{ "D" = "14.14"; "E" = "  +4.20%  " }
    24 = @{ "D" = "24.609.49"; "E" = "  +1.63%  " }
    25 = @{ "D" = "3.046"; "E" = "  +13.11%  " }
    26 = @{ "D" = "2.307"; "E" = "  -0.21%  " }
    27 = @{ "D" = "22.31"; "E" = "  +2.76%  " }
    28 = @{ "D" = "159.34"; "E" = "  +0.76%  " }
    29 = @{ "D" = "5.169"; "E" = "  -1.04%  " }
    30 = @{ "D" = "133.38"; "E" = "  +1.48%  " }
    31 = @{ "D" = "7.603"; "E" = "  +34.28%  " }
    32 = @{ "D" = "1.884.28"; "E" = "  +2.24%  " }
    33 = @{ "D" = "1.094"; "E" = "  -2.95%  " }
    34 = @{ "D" = "0.08621"; "E" = "  +1.55%  " }
    35 = @{ "D" = "7.338"; "E" = "  +21.81%  " }
    36 = @{ "B" = "WEMIXTOKEN"; "C" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; "D" = "1.962"; "E" = "  +10.32%  " }
    37 = @{ "B" = "FraxShare"; "C" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; "D" = "11.05"; "E" = "  +7.61%  " }
    38 = @{ "D" = "0.2719"; "E" = "  +4.42%  " }
    39 = @{ "D" = "14.75"; "E" = "  -0.11%  " }
    40 = @{ "E" = "  +10.81%  " }
    41 = @{ "D" = "0.08975"; "E" = "  +2.37%  " }
    42 = @{ "D" = "1.470"; "E" = "  +3.04%  " }
    43 = @{ "D" = "0.7647"; "E" = "  +5.00%  " }
    44 = @{ "D" = "0.7174"; "E" = "  +4.06%  " }
    45 = @{ "D" = "15.42"; "E" = "  +3.62%  " }
    46 = @{ "D" = "2.442"; "E" = "  +4.01%  " }
    47 = @{ "D" = "4.172"; "E" = "  +2.46%  " }
    48 = @{ "D" = "0.9979"; "E" = "  -0.87%  " }
    49 = @{ "D" = "140.63"; "E" = "  +1.60%  " }
    50 = @{ "B" = "Flow"; "C" = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"; "D" = "1.320"; "E" = "  +17.81%  " }
    51 = @{ "B" = "BabyDogeCoin"; "C" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; "D" = "0.00000000385"; "E" = "  +4.40%  " }
}

foreach ($row in $rowUpdates.Keys) {
    $cells = $rowUpdates[$row]
    foreach ($col in $cells.Keys) {
        $cellRef = "$col$row"
        $value = $cells[$col]

        # Column D holds price strings such as "0.9998" or "313.00" that Excel would
        # otherwise auto-coerce into numbers; force the cell to Text first so the
        # literal string (matching the site formatting) is preserved, just like the
        # multi-dot price strings (e.g. "24.625.33") that survive unchanged.
        if ($col -eq "D" -and $value -match "^-?\d+(\.\d+)?$") {
            $ws.Range($cellRef).NumberFormat = "@"
        }

        $ws.Range($cellRef).Value = $value
    }
}
